$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 267.875
$ws.Range("I2").Value = 234.71428
$ws.Range("K2").Value = 234.71428
$ws.Range("M2").Value = -121.71428
$ws.Range("H33").Value = 929.5833
$ws.Range("I33").Value = 220.25
$ws.Range("K33").Value = 220.25
$ws.Range("M33").Value = 8.75
$ws.Range("H70").Value = 3701
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 3551.5
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 10654.5
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -11194.5
$ws.Range("H73").Value = 3701
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 3551.5
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 10654.5
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -12526.5
$ws.Range("H106").Value = 725
$ws.Range("I106").Value = 725
$ws.Range("K106").Value = 725
$ws.Range("M106").Value = -94
$ws.Range("H107").Value = 649.8823
$ws.Range("I107").Value = 522.8461
$ws.Range("K107").Value = 522.8461
$ws.Range("M107").Value = 1397.1539
$ws.Range("H115").Value = 1379.25
$ws.Range("I115").Value = 1379.25
$ws.Range("K115").Value = 4137.75
$ws.Range("M115").Value = -2570.75

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1039.2727
$ws.Range("I2").Value = 1013.2
$ws.Range("K2").Value = 1013.2
$ws.Range("M2").Value = -900.2
$ws.Range("H24").Value = 78601.5
$ws.Range("J24").Value = 78601.5
$ws.Range("L24").Value = 78601.5
$ws.Range("N24").Value = -79349.5
$ws.Range("H25").Value = 805.3333
$ws.Range("I25").Value = 362.8
$ws.Range("K25").Value = 362.8
$ws.Range("M25").Value = 39.19999999999999
$ws.Range("H95").Value = 21475
$ws.Range("J95").Value = 21475
$ws.Range("L95").Value = 21475
$ws.Range("N95").Value = -26967
$ws.Range("H100").Value = 78601.5
$ws.Range("J100").Value = 78601.5
$ws.Range("L100").Value = 78601.5
$ws.Range("N100").Value = -80765.5
$ws.Range("H110").Value = 2968.9412
$ws.Range("I110").Value = 2037.9231
$ws.Range("J110").Value = 5994.75
$ws.Range("K110").Value = 2037.9231
$ws.Range("L110").Value = 5994.75
$ws.Range("M110").Value = 7.076900000000023
$ws.Range("N110").Value = -10084.75
$ws.Range("H116").Value = 1039.2727
$ws.Range("I116").Value = 1013.2
$ws.Range("K116").Value = 1013.2
$ws.Range("M116").Value = 1280.8

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1039.2727
$ws.Range("I3").Value = 1013.2
$ws.Range("K3").Value = 1013.2
$ws.Range("M3").Value = -899.2
$ws.Range("H86").Value = 1598.5
$ws.Range("I86").Value = 1598.5
$ws.Range("K86").Value = 1598.5
$ws.Range("M86").Value = -475.5
$ws.Range("H89").Value = 1598.5
$ws.Range("I89").Value = 1598.5
$ws.Range("K89").Value = 7992.5
$ws.Range("M89").Value = -2376.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 690.7646999999999
$ws.Range("I107").Value = 438.85715
$ws.Range("K107").Value = 438.85715
$ws.Range("M107").Value = 1481.14285

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 60
$ws.Range("I13").Value = 60
$ws.Range("K13").Value = 180
$ws.Range("M13").Value = -12

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 32521.076
$ws.Range("J24").Value = 32521.076
$ws.Range("L24").Value = 32521.076
$ws.Range("N24").Value = -32867.076
$ws.Range("H107").Value = 612.2857
$ws.Range("I107").Value = 612.2857
$ws.Range("K107").Value = 612.2857
$ws.Range("M107").Value = 1307.7143
$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11101.4
$ws.Range("I40").Value = 10472.235
$ws.Range("J40").Value = 14666.667
$ws.Range("K40").Value = 10472.235
$ws.Range("L40").Value = 14666.667
$ws.Range("M40").Value = -10336.235
$ws.Range("N40").Value = -14938.667
$ws.Range("H92").Value = 22111
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 969
$ws.Range("I93").Value = 797.1667
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 797.1667
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 450.8333
$ws.Range("N93").Value = -4496

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7899.8
$ws.Range("I62").Value = 7899.8
$ws.Range("K62").Value = 7899.8
$ws.Range("M62").Value = -7275.8
$ws.Range("H65").Value = 7899.8
$ws.Range("I65").Value = 7899.8
$ws.Range("K65").Value = 39499
$ws.Range("M65").Value = -36379
$ws.Range("H81").Value = 834
$ws.Range("I81").Value = 750
$ws.Range("K81").Value = 1500
$ws.Range("M81").Value = -439
$ws.Range("H84").Value = 834
$ws.Range("I84").Value = 750
$ws.Range("K84").Value = 7500
$ws.Range("M84").Value = -2196
$ws.Range("H87").Value = 70325
$ws.Range("J87").Value = 70325
$ws.Range("L87").Value = 70325
$ws.Range("N87").Value = -72821
$ws.Range("H90").Value = 70325
$ws.Range("J90").Value = 70325
$ws.Range("L90").Value = 210975
$ws.Range("N90").Value = -223455
$ws.Range("H107").Value = 892.36365
$ws.Range("I107").Value = 827
$ws.Range("K107").Value = 2481
$ws.Range("M107").Value = -561
$ws.Range("H122").Value = 893
$ws.Range("I122").Value = 893
$ws.Range("K122").Value = 2679
$ws.Range("M122").Value = -229
$ws.Range("H126").Value = 1988.9166
$ws.Range("I126").Value = 1586.8
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 4760.4
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -2290.4
$ws.Range("N126").Value = -16938.5
$ws.Range("H132").Value = 4938.8
$ws.Range("I132").Value = 4933
$ws.Range("K132").Value = 14799
$ws.Range("M132").Value = -12269
$ws.Range("H136").Value = 22528.44
$ws.Range("I136").Value = 23835.305
$ws.Range("J136").Value = 7499.5
$ws.Range("K136").Value = 71505.91500000001
$ws.Range("L136").Value = 22498.5
$ws.Range("M136").Value = -68955.91500000001
$ws.Range("N136").Value = -27598.5
